$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($row, $col, $value) {
    $c = $ws.Cells.Item($row, $col)
    $c.NumberFormat = "@"
    $c.Value = $value
    $c.Style = "Normal"
}

Set-TextCell 2 4 '29.562.15'
Set-TextCell 2 5 '  +2.06%  '
Set-TextCell 3 4 '1.988.04'
Set-TextCell 3 5 '  +5.28%  '
Set-TextCell 4 4 '0.9984'
Set-TextCell 4 5 '  -0.33%  '
Set-TextCell 5 4 '327.37'
Set-TextCell 5 5 '  +0.53%  '
Set-TextCell 6 4 '0.9980'
Set-TextCell 6 5 '  -0.34%  '
Set-TextCell 7 4 '0.4647'
Set-TextCell 7 5 '  +1.62%  '
Set-TextCell 8 4 '0.3948'
Set-TextCell 8 5 '  +1.09%  '
Set-TextCell 9 2 'Dogecoin'
Set-TextCell 9 3 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
Set-TextCell 9 4 '0.07936'
Set-TextCell 9 5 '  +0.92%  '
Set-TextCell 10 2 'Polygon'
Set-TextCell 10 3 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
Set-TextCell 10 4 '1.003'
Set-TextCell 10 5 '  +1.60%  '
Set-TextCell 11 2 'Solana'
Set-TextCell 11 3 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
Set-TextCell 11 4 '22.51'
Set-TextCell 11 5 '  +2.61%  '
Set-TextCell 12 2 'WrappedEther'
Set-TextCell 12 3 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
Set-TextCell 12 4 '1.979.69'
Set-TextCell 12 5 '  +7.09%  '
Set-TextCell 13 2 'Chainlink'
Set-TextCell 13 3 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
Set-TextCell 13 4 '7.213'
Set-TextCell 13 5 '  +2.49%  '
Set-TextCell 14 2 'Polkadot'
Set-TextCell 14 3 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
Set-TextCell 14 4 '5.857'
Set-TextCell 14 5 '  +2.78%  '
Set-TextCell 15 2 'TRON'
Set-TextCell 15 3 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
Set-TextCell 15 4 '0.07085'
Set-TextCell 15 5 '  +2.13%  '
Set-TextCell 16 2 'Litecoin'
Set-TextCell 16 3 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
Set-TextCell 16 4 '88.76'
Set-TextCell 16 5 '  +0.72%  '
Set-TextCell 17 2 'BinanceUSD'
Set-TextCell 17 3 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
Set-TextCell 17 4 '0.9995'
Set-TextCell 17 5 '  -0.29%  '
Set-TextCell 18 2 'ShibaInu'
Set-TextCell 18 3 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
Set-TextCell 18 4 '0.000009981'
Set-TextCell 18 5 '  -0.06%  '
Set-TextCell 19 2 'Avalanche'
Set-TextCell 19 3 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
Set-TextCell 19 4 '17.17'
Set-TextCell 19 5 '  +0.79%  '
Set-TextCell 20 2 'Dai'
Set-TextCell 20 3 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
Set-TextCell 20 4 '0.9980'
Set-TextCell 20 5 '  -0.38%  '
Set-TextCell 21 2 'WrappedBTC'
Set-TextCell 21 3 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
Set-TextCell 21 4 '29.600.57'
Set-TextCell 21 5 '  +2.22%  '
Set-TextCell 22 2 'Uniswap'
Set-TextCell 22 3 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
Set-TextCell 22 4 '5.541'
Set-TextCell 22 5 '  +4.73%  '
Set-TextCell 23 2 'Cosmos'
Set-TextCell 23 3 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
Set-TextCell 23 4 '11.26'
Set-TextCell 23 5 '  +2.70%  '
Set-TextCell 24 2 'WrappedliquidstakedEther2.0'
Set-TextCell 24 3 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
Set-TextCell 24 4 '2.235.46'
Set-TextCell 24 5 '  +7.62%  '
Set-TextCell 25 2 'Toncoin'
Set-TextCell 25 3 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
Set-TextCell 25 4 '2.122'
Set-TextCell 25 5 '  +3.40%  '
Set-TextCell 26 2 'Monero'
Set-TextCell 26 3 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
Set-TextCell 26 4 '157.60'
Set-TextCell 26 5 '  +1.19%  '
Set-TextCell 27 2 'EthereumClassic'
Set-TextCell 27 3 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
Set-TextCell 27 4 '19.63'
Set-TextCell 27 5 '  +1.70%  '
Set-TextCell 28 2 'InternetComputer(DFINITY)'
Set-TextCell 28 3 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
Set-TextCell 28 4 '6.024'
Set-TextCell 28 5 '  +0.92%  '
Set-TextCell 29 2 'BitcoinCash'
Set-TextCell 29 3 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
Set-TextCell 29 4 '120.40'
Set-TextCell 29 5 '  +2.31%  '
Set-TextCell 30 2 'LidoDAOToken'
Set-TextCell 30 3 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
Set-TextCell 30 4 '1.927'
Set-TextCell 30 5 '  +0.08%  '
Set-TextCell 31 2 'Stellar'
Set-TextCell 31 3 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
Set-TextCell 31 4 '0.09421'
Set-TextCell 31 5 '  +1.16%  '
Set-TextCell 32 2 'ImmutableX'
Set-TextCell 32 3 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
Set-TextCell 32 4 '0.8931'
Set-TextCell 32 5 '  -1.55%  '
Set-TextCell 33 2 'PEPE'
Set-TextCell 33 3 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
Set-TextCell 33 4 '0.000004198'
Set-TextCell 33 5 '  +154.68%  '
Set-TextCell 34 2 'Filecoin'
Set-TextCell 34 3 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
Set-TextCell 34 4 '5.271'
Set-TextCell 34 5 '  -0.37%  '
Set-TextCell 35 2 'ARBITRUM'
Set-TextCell 35 3 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
Set-TextCell 35 4 '1.349'
Set-TextCell 35 5 '  +1.37%  '
Set-TextCell 36 2 'HuobiToken'
Set-TextCell 36 3 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
Set-TextCell 36 4 '3.164'
Set-TextCell 36 5 '  -3.00%  '
Set-TextCell 37 2 'Hedera'
Set-TextCell 37 3 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
Set-TextCell 37 4 '0.05820'
Set-TextCell 37 5 '  +0.94%  '
Set-TextCell 38 2 'TrustWalletToken'
Set-TextCell 38 3 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
Set-TextCell 38 4 '1.178'
Set-TextCell 38 5 '  -2.10%  '
Set-TextCell 39 2 'VeChain'
Set-TextCell 39 3 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
Set-TextCell 39 4 '0.02128'
Set-TextCell 39 5 '  +2.68%  '
Set-TextCell 40 2 'FraxShare'
Set-TextCell 40 3 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
Set-TextCell 40 4 '7.923'
Set-TextCell 40 5 '  +3.66%  '
Set-TextCell 41 2 'Frax'
Set-TextCell 41 3 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
Set-TextCell 41 4 '0.9974'
Set-TextCell 41 5 '  -0.40%  '
Set-TextCell 42 2 'TheSandbox'
Set-TextCell 42 3 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
Set-TextCell 42 4 '0.5765'
Set-TextCell 42 5 '  +1.52%  '
Set-TextCell 43 2 'Algorand'
Set-TextCell 43 3 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
Set-TextCell 43 4 '0.1823'
Set-TextCell 43 5 '  +2.90%  '
Set-TextCell 44 2 'Aptos'
Set-TextCell 44 3 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
Set-TextCell 44 4 '9.823'
Set-TextCell 44 5 '  +0.79%  '
Set-TextCell 45 2 'EnergySwap'
Set-TextCell 45 3 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
Set-TextCell 45 4 '12.09'
Set-TextCell 45 5 '  +1.14%  '
Set-TextCell 46 2 'Decentraland'
Set-TextCell 46 3 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
Set-TextCell 46 4 '0.5387'
Set-TextCell 46 5 '  +0.41%  '
Set-TextCell 47 2 'RenderToken'
Set-TextCell 47 3 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
Set-TextCell 47 4 '2.181'
Set-TextCell 47 5 '  -4.08%  '
Set-TextCell 48 2 'MXToken'
Set-TextCell 48 3 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
Set-TextCell 48 4 '2.636'
Set-TextCell 48 5 '  +4.75%  '
Set-TextCell 49 2 'NEARProtocol'
Set-TextCell 49 3 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
Set-TextCell 49 4 '1.877'
Set-TextCell 49 5 '  +1.12%  '
Set-TextCell 50 4 '0.06977'
Set-TextCell 50 5 '  -0.84%  '
Set-TextCell 51 2 'Quant'
Set-TextCell 51 3 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
Set-TextCell 51 4 '114.32'
Set-TextCell 51 5 '  +1.55%  '
